# Apply "Add data for 2022-07-27" update:
#  - rename sheet/title date from 07-18 to 07-19
#  - update the "July (through 07-18)" label to "July (through 07-19)"
#  - bump June 2022 (I7) by 1
#  - update July row (row 8) figures for columns C..I
#  - update Total row (row 9) figures for columns C..I accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet (tab) name
$ws.Name = "Through 2022-07-19"

# Shared-string label for the July row
$ws.Range("A8").Value = "July (through 07-19)"

# June 2022 figure
$ws.Range("I7").Value = 143

# July row (2016 .. 2022)
$ws.Range("C8").Value = 38
$ws.Range("D8").Value = 39
$ws.Range("E8").Value = 47
$ws.Range("F8").Value = 29
$ws.Range("G8").Value = 79
$ws.Range("H8").Value = 95
$ws.Range("I8").Value = 105

# Total row (2016 .. 2022)
$ws.Range("C9").Value = 286
$ws.Range("D9").Value = 429
$ws.Range("E9").Value = 400
$ws.Range("F9").Value = 280
$ws.Range("G9").Value = 551
$ws.Range("H9").Value = 855
$ws.Range("I9").Value = 911
